# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# sheets to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 305
$ws1.Range("F10").Value = 4652
$ws1.Range("F19").Value = 95
$ws1.Range("F20").Value = 3540
$ws1.Range("F22").Value = 567
$ws1.Range("F29").Value = 73
$ws1.Range("F32").Value = 752
$ws1.Range("F33").Value = 2184
$ws1.Range("F34").Value = 403

# --- Sheet "全部类型" (all types, combined exhibitions + shows) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 305
$ws4.Range("F10").Value = 4652
$ws4.Range("F19").Value = 95
$ws4.Range("F20").Value = 3540
$ws4.Range("F22").Value = 567
$ws4.Range("F29").Value = 73
$ws4.Range("F33").Value = 752
$ws4.Range("F34").Value = 2184
$ws4.Range("F35").Value = 403
